$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 19451.102
$ws.Range("A12").Value = 18061.421
$ws.Range("A13").Value = 18201.4465
$ws.Range("A14").Value = 19619.5785
$ws.Range("A15").Value = 20554.4725
$ws.Range("A16").Value = 21841.412
$ws.Range("A17").Value = 22296.222
$ws.Range("A18").Value = 21879.0265
$ws.Range("A19").Value = 24784.9395
$ws.Range("A20").Value = 27200.6845
$ws.Range("A21").Value = 26829.0295
$ws.Range("A22").Value = 29178.267
$ws.Range("A23").Value = 31911.4905
